$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# ---------------------------------------------------------------------------
# 1. Insert a new row at row 9 (pushes the existing rows 9-18 down to 10-19)
#    and copy the formatting from row 8 (the "No Project" NP01 row) onto it,
#    since the new row is another "No Project" entry with identical styling.
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()

$ws.Range("A8:J8").Copy()
$ws.Range("A9:J9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Fill in the values for the newly inserted row 9:
#    2035_TM152_NGF_NP02 (No Project)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "NextGenFwys"
$ws.Range("B9").Value = 2035
$ws.Range("C9").Value = "2035_TM152_NGF_NP02"
$ws.Range("D9").Value = "NGF"
$ws.Range("E9").Value = "No Project"
$ws.Range("F9").Value = "`"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION`""
$ws.Range("G9").Value = "run182"
$ws.Range("H9").Value = "current"
$ws.Range("I9").Value = "NGF_Networks_NoProject_03"
$ws.Range("J9").Value = "https://app.asana.com/0/0/1202503525669953/f"

# ---------------------------------------------------------------------------
# 3. The original NP01 row (row 8) no longer shows a "current" status.
# ---------------------------------------------------------------------------
$ws.Range("H8").ClearContents()

# ---------------------------------------------------------------------------
# 4. The row that used to hold "...BlueprintSegmentedTest" (originally row 18,
#    now row 19 after the insert) is renamed to the new BPALT segmented run.
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = "2035_TM152_NGF_NP02_BPALTsegmented_00"

# ---------------------------------------------------------------------------
# 5. Re-point the two hyperlinks so they stay attached to the same rows'
#    content after the row shift (Blueprint_00 row moved from 17 -> 18;
#    the NP01 row stayed at J8). Recreate them with their original targets.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J18"), "https://app.asana.com/0/0/1202521542566668/f")
$ws.Range("J18").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("J8"), "https://app.asana.com/0/search?q=2035_TM152_NGF_NoProject_01&child=1201295328698176")
$ws.Range("J8").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 6. Restore the user's view state: selected cell D24 on the frozen pane.
# ---------------------------------------------------------------------------
$ws.Range("D24").Select()
